$wb = $excel.ActiveWorkbook

$wsInput    = $wb.Worksheets.Item("NewLoanInput")
$wsSummary  = $wb.Worksheets.Item("Summary")
$wsSchedule = $wb.Worksheets.Item("Repayment Schedule")

# --- Repayment Schedule: row 2 data tweaks ---------------------------------
# Pull formatting (style) from A2 (style index 8) onto the newly-populated
# B2 / F2 / O2 cells so they match the rest of the row instead of picking up
# the default style.
$wsSchedule.Cells.Item(2, 1).Copy()
$wsSchedule.Cells.Item(2, 2).PasteSpecial(-4122)   # B2 -> xlPasteFormats
$wsSchedule.Cells.Item(2, 6).PasteSpecial(-4122)   # F2 -> xlPasteFormats
$wsSchedule.Cells.Item(2, 15).PasteSpecial(-4122)  # O2 -> xlPasteFormats
$excel.CutCopyMode = $false

# H2: 0 -> blank, J2: 0 -> blank
$wsSchedule.Cells.Item(2, 8).ClearContents()
$wsSchedule.Cells.Item(2, 10).ClearContents()

# I2: blank -> 0, L2: blank -> 0
$wsSchedule.Cells.Item(2, 9).Value = 0
$wsSchedule.Cells.Item(2, 12).Value = 0

# --- Selections on each sheet ----------------------------------------------
# Update the Summary sheet's remembered selection without leaving it as the
# active tab.
$wsSummary.Range("D4").Select()

# Finally select on the Repayment Schedule sheet and leave it active - this
# is the sheet/tab that ends up selected when the workbook is reopened.
$wsSchedule.Range("D9").Select()
